$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - RandomForestRegressor (name unchanged), update values
$ws.Range("B3").Value = 0.9888885749447098
$ws.Range("C3").Value = 0.9877967006734268
$ws.Range("D3").Value = 0.834683084050041

# Row 4 - GradientBoostingRegressor -> DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.9837434359114132
$ws.Range("C4").Value = 0.9823441979829119
$ws.Range("D4").Value = 0.8475622989177116

# Row 5 - AdaBoostRegressor -> MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.8583893558900607
$ws.Range("C5").Value = 0.8479154043561267
$ws.Range("D5").Value = 0.7727108216995903
